$wb = $excel.ActiveWorkbook

# --- Sheet: Matriz_Económico (pairwise comparison matrix update) ---
$wsMatriz = $wb.Worksheets.Item("Matriz_Económico")
$wsMatriz.Range("D2").Value = 0.3333333333333333
$wsMatriz.Range("F2").Value = 0.3333333333333333
$wsMatriz.Range("G2").Value = 0.3333333333333333
$wsMatriz.Range("M2").Value = 0.3333333333333333
$wsMatriz.Range("N2").Value = 7
$wsMatriz.Range("D3").Value = 0.3333333333333333
$wsMatriz.Range("F3").Value = 0.3333333333333333
$wsMatriz.Range("G3").Value = 0.3333333333333333
$wsMatriz.Range("M3").Value = 0.3333333333333333
$wsMatriz.Range("N3").Value = 7
$wsMatriz.Range("B4").Value = 3
$wsMatriz.Range("C4").Value = 3
$wsMatriz.Range("E4").Value = 3
$wsMatriz.Range("J4").Value = 3
$wsMatriz.Range("K4").Value = 3
$wsMatriz.Range("D5").Value = 0.3333333333333333
$wsMatriz.Range("F5").Value = 0.3333333333333333
$wsMatriz.Range("G5").Value = 0.3333333333333333
$wsMatriz.Range("M5").Value = 0.3333333333333333
$wsMatriz.Range("N5").Value = 7
$wsMatriz.Range("B6").Value = 3
$wsMatriz.Range("C6").Value = 3
$wsMatriz.Range("E6").Value = 3
$wsMatriz.Range("J6").Value = 3
$wsMatriz.Range("K6").Value = 3
$wsMatriz.Range("B7").Value = 3
$wsMatriz.Range("C7").Value = 3
$wsMatriz.Range("E7").Value = 3
$wsMatriz.Range("J7").Value = 3
$wsMatriz.Range("K7").Value = 3
$wsMatriz.Range("D10").Value = 0.3333333333333333
$wsMatriz.Range("F10").Value = 0.3333333333333333
$wsMatriz.Range("G10").Value = 0.3333333333333333
$wsMatriz.Range("M10").Value = 0.3333333333333333
$wsMatriz.Range("N10").Value = 7
$wsMatriz.Range("D11").Value = 0.3333333333333333
$wsMatriz.Range("F11").Value = 0.3333333333333333
$wsMatriz.Range("G11").Value = 0.3333333333333333
$wsMatriz.Range("M11").Value = 0.3333333333333333
$wsMatriz.Range("B13").Value = 3
$wsMatriz.Range("C13").Value = 3
$wsMatriz.Range("E13").Value = 3
$wsMatriz.Range("J13").Value = 3
$wsMatriz.Range("K13").Value = 3
$wsMatriz.Range("B14").Value = 0.1428571428571428
$wsMatriz.Range("C14").Value = 0.1428571428571428
$wsMatriz.Range("E14").Value = 0.1428571428571428
$wsMatriz.Range("J14").Value = 0.1428571428571428

# --- Sheet: Pesos_Locales_Económico (recomputed local weights) ---
$wsPesos = $wb.Worksheets.Item("Pesos_Locales_Económico")
$wsPesos.Range("B2").Value = 0.06859393436079969
$wsPesos.Range("B3").Value = 0.06859393436079969
$wsPesos.Range("B4").Value = 0.1451001396860344
$wsPesos.Range("B5").Value = 0.06859393436079965
$wsPesos.Range("B6").Value = 0.1451001396860344
$wsPesos.Range("B7").Value = 0.1451001396860344
$wsPesos.Range("B8").Value = 0.01428744057464852
$wsPesos.Range("B9").Value = 0.01428744057464852
$wsPesos.Range("B10").Value = 0.06859393436079965
$wsPesos.Range("B11").Value = 0.0697866565851307
$wsPesos.Range("B12").Value = 0.008867907233478526
$wsPesos.Range("B13").Value = 0.1451001396860344
$wsPesos.Range("B14").Value = 0.009419377695460395
$wsPesos.Range("B15").Value = 0.01428744057464852
$wsPesos.Range("B16").Value = 0.01428744057464852

# --- Sheet: Resultados (recomputed global weights) ---
$wsResultados = $wb.Worksheets.Item("Resultados")
$wsResultados.Range("B2").Value = 0.05005747602085248
$wsResultados.Range("B3").Value = 0.02973933136067251
$wsResultados.Range("B4").Value = 0.0505184751493646
$wsResultados.Range("B5").Value = 0.1345809260762366
$wsResultados.Range("B6").Value = 0.03051182994464846
$wsResultados.Range("B7").Value = 0.03177932799136213
$wsResultados.Range("B8").Value = 0.0913964377575122
$wsResultados.Range("B9").Value = 0.01865873016652775
$wsResultados.Range("B10").Value = 0.0536879931345333
$wsResultados.Range("B11").Value = 0.1443287918427735
$wsResultados.Range("B12").Value = 0.1816769274842323
$wsResultados.Range("B13").Value = 0.03045392859632574
$wsResultados.Range("B14").Value = 0.05305423408053035
$wsResultados.Range("B15").Value = 0.0532741114550904
$wsResultados.Range("B16").Value = 0.04628147893933764

# --- Sheet: Ranking_Alternativas (recomputed ranking + reordered names) ---
$wsRanking = $wb.Worksheets.Item("Ranking_Alternativas")
$wsRanking.Range("B2").Value = 0.1816769274842323
$wsRanking.Range("B3").Value = 0.1443287918427735
$wsRanking.Range("B4").Value = 0.1345809260762366
$wsRanking.Range("B5").Value = 0.0913964377575122
$wsRanking.Range("B6").Value = 0.0536879931345333
$wsRanking.Range("B7").Value = 0.0532741114550904
$wsRanking.Range("B8").Value = 0.05305423408053035
$wsRanking.Range("B9").Value = 0.0505184751493646
$wsRanking.Range("B10").Value = 0.05005747602085248
$wsRanking.Range("B11").Value = 0.04628147893933764
$wsRanking.Range("B12").Value = 0.03177932799136213
$wsRanking.Range("B13").Value = 0.03051182994464846
$wsRanking.Range("B14").Value = 0.03045392859632574
$wsRanking.Range("B15").Value = 0.02973933136067251
$wsRanking.Range("B16").Value = 0.01865873016652775
$wsRanking.Range("A6").Value = "Placeres"
$wsRanking.Range("A7").Value = "Reina Isabel 2"
$wsRanking.Range("A8").Value = "Quebrada Verde"
